$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.986.27"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.417.11"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.10"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.10"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +4.71%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.848.91"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.931.82"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.431.61"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.76"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.06"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.61"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.56"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "328.69"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.406"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.42"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0517"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.577"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.404"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.05"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
